$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume (E) columns hold text values (e.g. "65.443.61", "0.997",
# "  +1.77%  ") in the source workbook. Temporarily mark the range as Text so that
# Excel does not auto-convert numeric-looking strings into real numbers, then restore
# the original (default/"Normal") cell style once the text has been written so no
# stray formatting is introduced.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '65.443.61'
$ws.Range("E2").Value = '  +1.77%  '
$ws.Range("D3").Value = '3.440.39'
$ws.Range("E3").Value = '  +3.71%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '550.91'
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("D6").Value = '179.88'
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("D7").Value = '0.642'
$ws.Range("E7").Value = '  +8.73%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +3.00%  '
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  +10.12%  '
$ws.Range("D11").Value = '53.54'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '0.0000271'
$ws.Range("E12").Value = '  +4.47%  '
$ws.Range("D13").Value = '9.16'
$ws.Range("E13").Value = '  +2.26%  '
$ws.Range("D14").Value = '3.988.74'
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = '0.121'
$ws.Range("E15").Value = '  +2.72%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.436.81'
$ws.Range("E16").Value = '  +3.72%  '
$ws.Range("D17").Value = '18.28'
$ws.Range("E17").Value = '  +4.61%  '
$ws.Range("D18").Value = '65.379.96'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").Value = '11.79'
$ws.Range("E19").Value = '  +5.21%  '
$ws.Range("E20").Value = '  +2.50%  '
$ws.Range("D21").Value = '415.86'
$ws.Range("E21").Value = '  +8.48%  '
$ws.Range("E22").Value = '  +7.87%  '
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '84.89'
$ws.Range("E24").Value = '  +3.49%  '
$ws.Range("E25").Value = '  -3.43%  '
$ws.Range("E26").Value = '  +4.17%  '
$ws.Range("D27").Value = '12.12'
$ws.Range("E27").Value = '  +7.42%  '
$ws.Range("D28").Value = '6.02'
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").Value = '8.88'
$ws.Range("E29").Value = '  +8.27%  '
$ws.Range("D30").Value = '29.74'
$ws.Range("E30").Value = '  +2.95%  '
$ws.Range("D31").Value = '6.50'
$ws.Range("E31").Value = '  -2.91%  '
$ws.Range("D32").Value = '611.63'
$ws.Range("E32").Value = '  -3.94%  '
$ws.Range("D33").Value = '11.64'
$ws.Range("E33").Value = '  +3.65%  '
$ws.Range("E34").Value = '  +3.25%  '
$ws.Range("D35").Value = '59.17'
$ws.Range("E35").Value = '  +3.30%  '
$ws.Range("E36").Value = '  +17.35%  '
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").Value = '37.22'
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("D39").Value = '0.0₃0778'
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").Value = '3.135.94'
$ws.Range("E41").Value = '  +5.70%  '
$ws.Range("D42").Value = '3.33'
$ws.Range("E42").Value = '  +2.40%  '
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  -3.59%  '
$ws.Range("D45").Value = '2.80'
$ws.Range("E45").Value = '  +5.37%  '
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").Value = '  +1.99%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.20'
$ws.Range("E47").Value = '  +2.02%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '2.72'
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("D49").Value = '0.131'
$ws.Range("E49").Value = '  +4.92%  '
$ws.Range("D50").Value = '138.18'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").Value = '8.37'
$ws.Range("E51").Value = '  +1.00%  '

$rng.Style = "Normal"
